$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The property-type header row (row 4) gains two new trailing columns:
# "Pattern" (M4) and "Pattern Type" (N4), matching the style already
# used by the neighboring "Unique" header cell (L4).
$ws.Range("L4").Copy()
$ws.Range("M4:N4").PasteSpecial(-4122)

$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

# Reflect the newly added cells in the active selection, like the
# original edit did (previously L4:L7 was selected).
$ws.Range("M4:N4").Select()
